$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Row 9: set D9 = "OK" and E9 = date 2010-04-14 (matches existing E4/E8 entries)
$ws.Range("D9").Value = "OK"
$ws.Range("E9").Value = (Get-Date -Year 2010 -Month 4 -Day 14)

# Update active selection to F9 as in the diff
$ws.Range("F9").Select()
